$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2487.3
$ws.Range("I98").Value = 2487.3
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2487.3
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -989.3000000000002

$ws.Range("H100").Value = 2600
$ws.Range("I100").Value = 2600
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2600
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2059
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 2487.3
$ws.Range("I122").Value = 2487.3
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7461.900000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5011.900000000001

$ws.Range("H133").Value = 99999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -110119

$ws.Range("H135").Value = 1172.25
$ws.Range("I135").Value = 522.7826
$ws.Range("J135").Value = 4159.8
$ws.Range("K135").Value = 4705.0434
$ws.Range("L135").Value = 37438.2
$ws.Range("M135").Value = -2170.0434

$ws.Range("H137").Value = 1373.6666
$ws.Range("I137").Value = 1309.6666
$ws.Range("J137").Value = 1693.6666
$ws.Range("K137").Value = 3928.9998
$ws.Range("L137").Value = 5080.9998
$ws.Range("M137").Value = -1378.9998

$ws.Range("H138").Value = 2746.689
$ws.Range("I138").Value = 1368.6111
$ws.Range("J138").Value = 3665.4075
$ws.Range("K138").Value = 4105.8333
$ws.Range("L138").Value = 10996.2225
$ws.Range("M138").Value = 1034.1667
$ws.Range("N138").Value = -21276.2225

$ws.Range("H141").Value = 3174.1
$ws.Range("I141").Value = 3135.6155
$ws.Range("J141").Value = 3424.25
$ws.Range("K141").Value = 9406.8465
$ws.Range("L141").Value = 10272.75
$ws.Range("M141").Value = -4226.8465
$ws.Range("N141").Value = -20632.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H38").Value = 741.8
$ws.Range("I38").Value = 741.8
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 741.8
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -274.8
$ws.Range("N38").ClearContents()

$ws.Range("H44").Value = 23540.834
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 23540.834
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 23540.834
$ws.Range("N44").Value = -24516.834

$ws.Range("H61").Value = 3610.8667
$ws.Range("I61").Value = 3252.353
$ws.Range("J61").Value = 4079.6924
$ws.Range("K61").Value = 3252.353
$ws.Range("L61").Value = 4079.6924
$ws.Range("M61").Value = -3040.353
$ws.Range("N61").Value = -4503.6924

$ws.Range("H122").Value = 3015.7144
$ws.Range("I122").Value = 2674
$ws.Range("J122").Value = 3471.3333
$ws.Range("K122").Value = 8022
$ws.Range("L122").Value = 10413.9999
$ws.Range("M122").Value = -5572
$ws.Range("N122").Value = -15313.9999

$ws.Range("H132").Value = 4101.0454
$ws.Range("I132").Value = 1975.1052
$ws.Range("J132").Value = 17565.334
$ws.Range("K132").Value = 5925.3156
$ws.Range("L132").Value = 52696.00199999999
$ws.Range("M132").Value = -3395.3156

$ws.Range("H136").Value = 3610.8667
$ws.Range("I136").Value = 3252.353
$ws.Range("J136").Value = 4079.6924
$ws.Range("K136").Value = 9757.059000000001
$ws.Range("L136").Value = 12239.0772
$ws.Range("M136").Value = -7207.059000000001
$ws.Range("N136").Value = -17339.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 22224684
$ws.Range("I134").Value = 2413.75
$ws.Range("J134").Value = 111113770
$ws.Range("K134").Value = 7241.25
$ws.Range("L134").Value = 333341310
$ws.Range("M134").Value = -4706.25
$ws.Range("N134").Value = -333346380

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 474.6
$ws.Range("I22").Value = 559.44446
$ws.Range("J22").Value = 347.33334
$ws.Range("K22").Value = 559.44446
$ws.Range("L22").Value = 347.33334
$ws.Range("M22").Value = -209.44446
$ws.Range("N22").Value = -1047.33334

$ws.Range("H31").Value = 2345
$ws.Range("I31").Value = 2204.4443
$ws.Range("J31").Value = 2766.6667
$ws.Range("K31").Value = 2204.4443
$ws.Range("L31").Value = 2766.6667
$ws.Range("M31").Value = -1909.4443
$ws.Range("N31").Value = -3356.6667

$ws.Range("H34").Value = 2345
$ws.Range("I34").Value = 2204.4443
$ws.Range("J34").Value = 2766.6667
$ws.Range("K34").Value = 2204.4443
$ws.Range("L34").Value = 2766.6667
$ws.Range("M34").Value = -2002.4443
$ws.Range("N34").Value = -3170.6667

$ws.Range("H58").Value = 2670.9375
$ws.Range("I58").Value = 2571.923
$ws.Range("J58").Value = 3100
$ws.Range("K58").Value = 2571.923
$ws.Range("L58").Value = 3100
$ws.Range("M58").Value = -2368.923

$ws.Range("H132").Value = 5041.8687
$ws.Range("I132").Value = 4850.6553
$ws.Range("J132").Value = 5658
$ws.Range("K132").Value = 14551.9659
$ws.Range("L132").Value = 16974
$ws.Range("M132").Value = -12021.9659
$ws.Range("N132").Value = -22034

$ws.Range("H134").Value = 5558189
$ws.Range("I134").Value = 2750.5
$ws.Range("J134").Value = 25002224
$ws.Range("K134").Value = 8251.5
$ws.Range("L134").Value = 75006672
$ws.Range("M134").Value = -5716.5
$ws.Range("N134").Value = -75011742

$ws.Range("H136").Value = 2670.9375
$ws.Range("I136").Value = 2571.923
$ws.Range("J136").Value = 3100
$ws.Range("K136").Value = 7715.768999999999
$ws.Range("L136").Value = 9300
$ws.Range("M136").Value = -5165.768999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4493.2
$ws.Range("I34").Value = 646.6667
$ws.Range("J34").Value = 6141.7144
$ws.Range("K34").Value = 1940.0001
$ws.Range("L34").Value = 18425.1432
$ws.Range("M34").Value = -1856.0001
$ws.Range("N34").Value = -18593.1432

$ws.Range("H55").Value = 423.5
$ws.Range("I55").Value = 308.2
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 924.5999999999999
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -747.5999999999999
$ws.Range("N55").Value = -3354

$ws.Range("H69").Value = 2757.6667
$ws.Range("I69").Value = 1506
$ws.Range("J69").Value = 3759
$ws.Range("K69").Value = 4518
$ws.Range("L69").Value = 11277
$ws.Range("M69").Value = -3707
$ws.Range("N69").Value = -12899

$ws.Range("H72").Value = 2757.6667
$ws.Range("I72").Value = 1506
$ws.Range("J72").Value = 3759
$ws.Range("K72").Value = 13554
$ws.Range("L72").Value = 33831
$ws.Range("M72").Value = -9498
$ws.Range("N72").Value = -41943

$ws.Range("H131").Value = 1551.5834
$ws.Range("I131").Value = 1269.8334
$ws.Range("J131").Value = 1833.3334
$ws.Range("K131").Value = 3809.5002
$ws.Range("L131").Value = 5500.0002
$ws.Range("M131").Value = 1230.4998
$ws.Range("N131").Value = -15580.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5066.8667
$ws.Range("I126").Value = 5850.3
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 17550.9
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -15080.9

$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2437.6572
$ws.Range("I46").Value = 1600
$ws.Range("J46").Value = 2545.742
$ws.Range("K46").Value = 1600
$ws.Range("L46").Value = 2545.742
$ws.Range("M46").Value = -1412

$ws.Range("H63").Value = 30027.889
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 30027.889
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 30027.889
$ws.Range("N63").Value = -31525.889

$ws.Range("H66").Value = 30027.889
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 30027.889
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 90083.667
$ws.Range("N66").Value = -97571.667

$ws.Range("H68").Value = 4436
$ws.Range("I68").Value = 4436
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4436
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -3687

$ws.Range("H71").Value = 4436
$ws.Range("I71").Value = 4436
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 22180
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -18436

$ws.Range("H80").Value = 20128
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 20128
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 20128
$ws.Range("N80").Value = -22374

$ws.Range("H83").Value = 20128
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 20128
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 60384
$ws.Range("N83").Value = -71616

$ws.Range("H132").Value = 2290.4
$ws.Range("I132").Value = 2290.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6871.200000000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4341.200000000001
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 55559260
$ws.Range("I136").Value = 3424.4546
$ws.Range("J136").Value = 142861280
$ws.Range("K136").Value = 10273.3638
$ws.Range("L136").Value = 428583840
$ws.Range("M136").Value = -7723.363799999999
$ws.Range("N136").Value = -428588940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 24599.7
$ws.Range("I2").Value = 9999.5
$ws.Range("J2").Value = 28249.75
$ws.Range("K2").Value = 9999.5
$ws.Range("L2").Value = 28249.75
$ws.Range("M2").Value = -9887.5
$ws.Range("N2").Value = -28473.75

$ws.Range("H4").Value = 6500
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 6500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 6500
$ws.Range("N4").Value = -6726

$ws.Range("H75").Value = 26500
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 26500
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 26500
$ws.Range("N75").Value = -28372

$ws.Range("H78").Value = 26500
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 26500
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 79500
$ws.Range("N78").Value = -88860

$ws.Range("H107").Value = 1058.2354
$ws.Range("I107").Value = 1066.1333
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 3198.3999
$ws.Range("L107").Value = 2997
$ws.Range("M107").Value = -1278.3999

$ws.Range("H122").Value = 4102.1
$ws.Range("I122").Value = 4113.5557
$ws.Range("J122").Value = 3999
$ws.Range("K122").Value = 12340.6671
$ws.Range("L122").Value = 11997
$ws.Range("M122").Value = -9890.667099999999

$ws.Range("H132").Value = 762.375
$ws.Range("I132").Value = 774.1429000000001
$ws.Range("J132").Value = 680
$ws.Range("K132").Value = 2322.4287
$ws.Range("L132").Value = 2040
$ws.Range("M132").Value = 207.5712999999996
$ws.Range("N132").Value = -7100
